# Update Riders (C) and Average (D) columns for Madigan bike hours
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$values = @(
    @{Row=2;  C=244; D=210.56},
    @{Row=3;  C=219; D=229.71},
    @{Row=4;  C=241; D=239.89},
    @{Row=5;  C=100; D=117.84},
    @{Row=6;  C=61;  D=99.47},
    @{Row=7;  C=212; D=228.94},
    @{Row=8;  C=216; D=211},
    @{Row=9;  C=239; D=212.05},
    @{Row=10; C=238; D=230.17},
    @{Row=11; C=219; D=238.79},
    @{Row=12; C=102; D=117.05},
    @{Row=13; C=71;  D=97.89},
    @{Row=14; C=267; D=231.06},
    @{Row=15; C=243; D=212.6},
    @{Row=16; C=206; D=211.75},
    @{Row=17; C=415; D=239.89},
    @{Row=18; C=213; D=237.5},
    @{Row=19; C=114; D=116.9},
    @{Row=20; C=74;  D=96.63},
    @{Row=21; C=277; D=233.47},
    @{Row=22; C=219; D=212.9},
    @{Row=23; C=215; D=211.9},
    @{Row=24; C=204; D=238.1},
    @{Row=25; C=305; D=240.71},
    @{Row=26; C=85;  D=115.45},
    @{Row=27; C=69;  D=95.25},
    @{Row=28; C=99;  D=226.75},
    @{Row=29; C=185; D=211.64},
    @{Row=30; C=221; D=212.32},
    @{Row=31; C=264; D=239.33},
    @{Row=32; C=260; D=241.59}
)

foreach ($row in $values) {
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $ws.Cells.Item($row.Row, 4).Value = $row.D
}
